$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.153.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.90"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5218"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07267"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8992"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08170"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.30%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.922.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.277"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008574"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.54"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.172.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.083"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.396"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Monero"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C24").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("B25").ClearFormats()
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C25").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.293"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.18"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.735"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.94"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.897"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.784"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09224"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05051"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7904"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.212"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("B35").ClearFormats()
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C35").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.973"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "MXToken"
$ws.Range("B36").ClearFormats()
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C36").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.410"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.582"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5668"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01986"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.990"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.550"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1517"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4861"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.618"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.09"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.18%  "
$ws.Range("E51").ClearFormats()

